$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")

# "add TAN to tag" - add Term Accession Number and Term Source REF values
# for the second tag ("Plant", column C) in the Tags section.
$ws.Range("C14").Value = "http://purl.obolibrary.org/obo/NCIT_C14258"
$ws.Range("C15").Value = "NCIT"
